# Auto-generated Excel COM-interop script applying the Unicorn_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for the
# affected leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 98.78570999999999
$ws.Range("I9").Value = 106.375
$ws.Range("J9").Value = 88.666664
$ws.Range("K9").Value = 106.375
$ws.Range("L9").Value = 88.666664
$ws.Range("M9").Value = 62.625
$ws.Range("N9").Value = -426.666664
$ws.Range("H15").Value = 105.7
$ws.Range("I15").Value = 105.7
$ws.Range("K15").Value = 317.1
$ws.Range("M15").Value = -148.1
$ws.Range("H140").Value = 46852.668
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 46852.668
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 46852.668
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -57212.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 2849
$ws.Range("I25").Value = 465.33334
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 465.33334
$ws.Range("L25").Value = 10000
$ws.Range("M25").Value = -63.33334000000002
$ws.Range("N25").Value = -10804
$ws.Range("H32").Value = 6166.25
$ws.Range("I32").Value = 3486.595
$ws.Range("J32").Value = 16246.857
$ws.Range("K32").Value = 3486.595
$ws.Range("L32").Value = 16246.857
$ws.Range("M32").Value = -3199.595
$ws.Range("N32").Value = -16820.857
$ws.Range("H61").Value = 2767.7896
$ws.Range("I61").Value = 1643.4
$ws.Range("J61").Value = 5413.4116
$ws.Range("K61").Value = 1643.4
$ws.Range("L61").Value = 5413.4116
$ws.Range("M61").Value = -1431.4
$ws.Range("N61").Value = -5837.4116
$ws.Range("H136").Value = 2767.7896
$ws.Range("I136").Value = 1643.4
$ws.Range("J136").Value = 5413.4116
$ws.Range("K136").Value = 4930.200000000001
$ws.Range("L136").Value = 16240.2348
$ws.Range("M136").Value = -2380.200000000001
$ws.Range("N136").Value = -21340.2348
$ws.Range("H139").Value = 22787.5
$ws.Range("J139").Value = 22787.5
$ws.Range("L139").Value = 22787.5
$ws.Range("N139").Value = -33067.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 5426.5
$ws.Range("I37").Value = 853
$ws.Range("J37").Value = 10000
$ws.Range("K37").Value = 853
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = -716
$ws.Range("N37").Value = -10274
$ws.Range("H138").Value = 55250
$ws.Range("J138").Value = 55250
$ws.Range("L138").Value = 55250
$ws.Range("N138").Value = -65530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 509334.34
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 509334.34
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 509334.34
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -509558.34
$ws.Range("H7").Value = 89.42856999999999
$ws.Range("I7").Value = 70.181816
$ws.Range("K7").Value = 70.181816
$ws.Range("M7").Value = 42.818184
$ws.Range("H22").Value = 642.7826
$ws.Range("I22").Value = 321.41666
$ws.Range("K22").Value = 321.41666
$ws.Range("M22").Value = 28.58334000000002
$ws.Range("H82").Value = 28720
$ws.Range("I82").Value = 9160
$ws.Range("J82").Value = 38500
$ws.Range("K82").Value = 9160
$ws.Range("L82").Value = 38500
$ws.Range("M82").Value = -8799
$ws.Range("N82").Value = -39222
$ws.Range("H85").Value = 28720
$ws.Range("I85").Value = 9160
$ws.Range("J85").Value = 38500
$ws.Range("K85").Value = 9160
$ws.Range("L85").Value = 38500
$ws.Range("M85").Value = -7912
$ws.Range("N85").Value = -40996
$ws.Range("H138").Value = 25285.715
$ws.Range("J138").Value = 25285.715
$ws.Range("L138").Value = 25285.715
$ws.Range("N138").Value = -35565.715
$ws.Range("H140").Value = 20300
$ws.Range("I140").Value = 9800
$ws.Range("J140").Value = 22400
$ws.Range("K140").Value = 9800
$ws.Range("L140").Value = 22400
$ws.Range("M140").Value = -4620
$ws.Range("N140").Value = -32760

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 907.5714
$ws.Range("I4").Value = 350
$ws.Range("J4").Value = 1000.5
$ws.Range("K4").Value = 1050
$ws.Range("L4").Value = 3001.5
$ws.Range("M4").Value = -938
$ws.Range("N4").Value = -3225.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 6500
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 12000
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = -888
$ws.Range("N5").Value = -12224
$ws.Range("H33").Value = 4990
$ws.Range("J33").Value = 4990
$ws.Range("L33").Value = 4990
$ws.Range("N33").Value = -5494
$ws.Range("H107").Value = 4436.2085
$ws.Range("I107").Value = 5143.25
$ws.Range("J107").Value = 901
$ws.Range("K107").Value = 5143.25
$ws.Range("L107").Value = 901
$ws.Range("M107").Value = -3223.25
$ws.Range("N107").Value = -4741
$ws.Range("H140").Value = 20534
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 20534
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 20534
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -30894

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 66673070
$ws.Range("J2").Value = 66673070
$ws.Range("L2").Value = 66673070
$ws.Range("N2").Value = -66673294
$ws.Range("H24").Value = 2900
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 2900
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 2900
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -3586
$ws.Range("H122").Value = 2829.2856
$ws.Range("I122").Value = 2669.6875
$ws.Range("J122").Value = 3340
$ws.Range("K122").Value = 8009.0625
$ws.Range("L122").Value = 10020
$ws.Range("M122").Value = -5559.0625
$ws.Range("N122").Value = -14920
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 8000
$ws.Range("I2").Value = 8000
$ws.Range("K2").Value = 8000
$ws.Range("M2").Value = -7888
$ws.Range("H39").Value = 9750
$ws.Range("J39").Value = 9750
$ws.Range("L39").Value = 9750
$ws.Range("N39").Value = -10576
$ws.Range("H41").Value = 3768.2727
$ws.Range("J41").Value = 3768.2727
$ws.Range("L41").Value = 3768.2727
$ws.Range("N41").Value = -4548.2727
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H43").Value = 12416.2
$ws.Range("I43").Value = 9750.5
$ws.Range("J43").Value = 14193.333
$ws.Range("K43").Value = 9750.5
$ws.Range("L43").Value = 14193.333
$ws.Range("M43").Value = -9601.5
$ws.Range("N43").Value = -14491.333
$ws.Range("H81").Value = 2238.9473
$ws.Range("I81").Value = 1407.5
$ws.Range("J81").Value = 2843.6365
$ws.Range("K81").Value = 2815
$ws.Range("L81").Value = 5687.273
$ws.Range("M81").Value = -1754
$ws.Range("N81").Value = -7809.273
$ws.Range("H84").Value = 2238.9473
$ws.Range("I84").Value = 1407.5
$ws.Range("J84").Value = 2843.6365
$ws.Range("K84").Value = 14075
$ws.Range("L84").Value = 28436.365
$ws.Range("M84").Value = -8771
$ws.Range("N84").Value = -39044.36500000001
$ws.Range("H141").Value = 36714.5
$ws.Range("J141").Value = 36714.5
$ws.Range("L141").Value = 36714.5
$ws.Range("N141").Value = -47074.5

Write-Output "Applied 196 cell updates across 8 sheets."
